{"js": "// Remove the standalone \"OBJECTIVE: \" heading paragraph that used to be the\n// first paragraph in the body, and relocate the \"_GoBack\" bookmark so it\n// wraps the start of the (new) first paragraph instead of sitting in the\n// middle of the big vitals paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The first paragraph in the document is the \"OBJECTIVE: \" heading \u2014 delete\n// it entirely (paragraph mark included).\nconst headingParagraph = paragraphs.items[0];\nheadingParagraph.delete();\nawait context.sync();\n\n// Re-query paragraphs now that the heading paragraph is gone so we get a\n// fresh reference to what is now the first paragraph in the body.\nconst remainingParagraphs = body.paragraphs;\nremainingParagraphs.load(\"items\");\nawait context.sync();\nconst firstParagraph = remainingParagraphs.items[0];\n\n// Drop the bookmark from its old spot (between \"Abdomen is soft and\n// non-tender to palpation. \" and \"Bowel sounds are normal.\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-insert the bookmark collapsed at the very start of the (new) first\n// paragraph.\nconst startRange = firstParagraph.getRange(\"Start\");\nstartRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Remove the standalone \"OBJECTIVE: \" heading paragraph that used to be the\n# first paragraph in the body, and relocate the \"_GoBack\" bookmark so it\n# wraps the start of the (new) first paragraph instead of sitting in the\n# middle of the big vitals paragraph.\n\n$d = $word.ActiveDocument\n\n# 1. Insert the new \"_GoBack\" bookmark collapsed at the very start of what is\n#    currently the second paragraph (\"Neuro: ...\"). Bookmarks.Add re-uses the\n#    name, so this automatically removes the old \"_GoBack\" bookmark that sat\n#    mid-document (between \"Abdomen is soft and non-tender to palpation. \"\n#    and \"Bowel sounds are normal.\").\n#    NOTE: we deliberately do this while paragraph 1 (\"OBJECTIVE: \") still\n#    precedes paragraph 2, so the insertion point is non-zero; deleting\n#    paragraph 1 afterwards shifts the (still-collapsed) bookmark down to the\n#    very start of the document.\n$secondParagraph = $d.Paragraphs.Item(2)\n$newBookmarkRange = $secondParagraph.Range.Duplicate\n$newBookmarkRange.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $newBookmarkRange)\n\n# 2. Delete the first paragraph (\"OBJECTIVE: \") entirely, paragraph mark\n#    included.\n$d.Paragraphs.Item(1).Range.Delete()\n"}
